$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to be treated as text so that numeric-looking
# strings (e.g. "1.001") are not silently converted to floating point
# numbers, matching the original inline-string cell content.
$ws.Range('D2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '22.477.11'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '1.572.20'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.13%  '
$ws.Range('D5').Value = '1.001'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').Value = '291.21'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').Value = '0.3753'
$ws.Range('E7').Value = '  -0.52%  '
$ws.Range('D8').Value = '49.85'
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('D9').Value = '0.3408'
$ws.Range('E9').Value = '  -0.53%  '
$ws.Range('D10').Value = '1.151'
$ws.Range('D11').Value = '0.07578'
$ws.Range('E11').Value = '  -1.46%  '
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').Value = '  -0.13%  '
$ws.Range('D13').Value = '21.42'
$ws.Range('E13').Value = '  +0.16%  '
$ws.Range('D14').Value = '6.019'
$ws.Range('E14').Value = '  +0.41%  '
$ws.Range('D15').Value = '6.961'
$ws.Range('E15').Value = '  +0.32%  '
$ws.Range('D16').Value = '1.572.76'
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').Value = '0.00001123'
$ws.Range('E17').Value = '  -1.59%  '
$ws.Range('D18').Value = '91.15'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('D19').Value = '0.06743'
$ws.Range('E19').Value = '  +0.06%  '
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = '6.273'
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('D22').Value = '16.42'
$ws.Range('E22').Value = '  -2.28%  '
$ws.Range('D23').Value = '12.21'
$ws.Range('E23').Value = '  +1.50%  '
$ws.Range('D24').Value = '22.478.68'
$ws.Range('E24').Value = '  +0.38%  '
$ws.Range('D25').Value = '2.331'
$ws.Range('E25').Value = '  -2.69%  '
$ws.Range('D26').Value = '2.602'
$ws.Range('E26').Value = '  -6.58%  '
$ws.Range('D27').Value = '20.17'
$ws.Range('E27').Value = '  -0.65%  '
$ws.Range('D28').Value = '148.38'
$ws.Range('E28').Value = '  +2.32%  '
$ws.Range('D29').Value = '4.998'
$ws.Range('E29').Value = '  -1.61%  '
$ws.Range('D30').Value = '126.17'
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('D31').Value = '1.749.15'
$ws.Range('E31').Value = '  +0.09%  '
$ws.Range('D32').Value = '1.051'
$ws.Range('E32').Value = '  +3.22%  '
$ws.Range('D33').Value = '6.164'
$ws.Range('D34').Value = '1.978'
$ws.Range('E34').Value = '  -2.46%  '
$ws.Range('D35').Value = '9.875'
$ws.Range('E35').Value = '  -2.16%  '
$ws.Range('D36').Value = '0.08464'
$ws.Range('E36').Value = '  -0.83%  '
$ws.Range('E37').Value = '  +6.17%  '
$ws.Range('D38').Value = '0.02475'
$ws.Range('E38').Value = '  -3.48%  '
$ws.Range('D39').Value = '0.2295'
$ws.Range('E39').Value = '  -1.39%  '
$ws.Range('D40').Value = '0.06567'
$ws.Range('E40').Value = '  -0.16%  '
$ws.Range('D41').Value = '5.492'
$ws.Range('E41').Value = '  -0.28%  '
$ws.Range('D42').Value = '11.43'
$ws.Range('E42').Value = '  -2.16%  '
$ws.Range('D43').Value = '0.6310'
$ws.Range('E43').Value = '  -2.42%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').Value = '14.11'
$ws.Range('E44').Value = '  +0.07%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').Value = '1.001'
$ws.Range('E45').Value = '  -0.05%  '
$ws.Range('D46').Value = '3.818'
$ws.Range('E46').Value = '  +0.86%  '
$ws.Range('D47').Value = '0.5895'
$ws.Range('E47').Value = '  -2.47%  '
$ws.Range('D48').Value = '2.103'
$ws.Range('E48').Value = '  -0.14%  '
$ws.Range('D49').Value = '129.91'
$ws.Range('E49').Value = '  +3.43%  '
$ws.Range('D50').Value = '1.233'
$ws.Range('E50').Value = '  -5.51%  '
$ws.Range('D51').Value = '0.07343'
$ws.Range('E51').Value = '  +0.03%  '

# Restore original (default) cell formatting/style so the saved
# worksheet matches the source style indices (no "s" attribute).
$ws.Range('D2:E51').ClearFormats()

